# Saldo_guide.xlsx update — refresh the daily client-balance export:
#  - rename the sheet to the new extraction timestamp
#  - bump every "Dt. Referencia" (column G) from 2024-09-02 to 2024-09-03
#  - refresh the few rows whose projected/expected balances changed between
#    the two extractions (rows 43, 138, 241)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet (also updates the <sheet name=.../> entry in workbook.xml)
$ws.Name = "IClientBalance-20240903-092801-"

# Column G holds the reference date for every data row (2..274), stored as
# the Excel serial date 45537 (2024-09-02) -> 45538 (2024-09-03).
for ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value2 = 45538
}

# Row 43: Vl. Projetado (D) and Saldo Previsto (E) both shift; Vl. Total (H)
# is their sum.
$ws.Cells.Item(43, 4).Value2 = 0
$ws.Cells.Item(43, 5).Value2 = 52783.43
$ws.Cells.Item(43, 8).Value2 = 52783.43

# Row 138: Saldo Previsto (E) updates; Vl. Total (H) mirrors it.
$ws.Cells.Item(138, 5).Value2 = 152809.07
$ws.Cells.Item(138, 8).Value2 = 152809.07

# Row 241: Saldo Previsto (E) updates; Vl. Total (H) mirrors it.
$ws.Cells.Item(241, 5).Value2 = 6524.84
$ws.Cells.Item(241, 8).Value2 = 6524.84
